$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 203.25
$ws.Range("I6").Value = 89.42856999999999
$ws.Range("K6").Value = 268.28571
$ws.Range("M6").Value = -156.28571

$ws.Range("H43").Value = 8000.5
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H70").Value = 2960.6924
$ws.Range("I70").Value = 4141.2856
$ws.Range("J70").Value = 1583.3334
$ws.Range("K70").Value = 12423.8568
$ws.Range("L70").Value = 4750.0002
$ws.Range("M70").Value = -12153.8568
$ws.Range("N70").Value = -5290.0002

$ws.Range("H73").Value = 2960.6924
$ws.Range("I73").Value = 4141.2856
$ws.Range("J73").Value = 1583.3334
$ws.Range("K73").Value = 12423.8568
$ws.Range("L73").Value = 4750.0002
$ws.Range("M73").Value = -11487.8568
$ws.Range("N73").Value = -6622.0002

$ws.Range("H74").Value = 3554.7
$ws.Range("I74").Value = 3554.7
$ws.Range("K74").Value = 3554.7
$ws.Range("M74").Value = -2618.7

$ws.Range("H76").Value = 28574320
$ws.Range("I76").Value = 50001836
$ws.Range("K76").Value = 50001836
$ws.Range("M76").Value = -50001521

$ws.Range("H77").Value = 3554.7
$ws.Range("I77").Value = 3554.7
$ws.Range("K77").Value = 17773.5
$ws.Range("M77").Value = -13093.5

$ws.Range("H79").Value = 28574320
$ws.Range("I79").Value = 50001836
$ws.Range("K79").Value = 50001836
$ws.Range("M79").Value = -50000744

$ws.Range("H98").Value = 1021.1111
$ws.Range("I98").Value = 1021.1111
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1021.1111
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 476.8889
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 2422.6667
$ws.Range("I100").Value = 2264.6924
$ws.Range("K100").Value = 2264.6924
$ws.Range("M100").Value = -1723.6924

$ws.Range("H122").Value = 1021.1111
$ws.Range("I122").Value = 1021.1111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3063.3333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -613.3332999999998
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2016.9
$ws.Range("I45").Value = 2016.9
$ws.Range("K45").Value = 2016.9
$ws.Range("M45").Value = -1639.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3625.9333
$ws.Range("I86").Value = 3953.077
$ws.Range("J86").Value = 1499.5
$ws.Range("K86").Value = 3953.077
$ws.Range("L86").Value = 1499.5
$ws.Range("M86").Value = -2830.077
$ws.Range("N86").Value = -3745.5

$ws.Range("H89").Value = 3625.9333
$ws.Range("I89").Value = 3953.077
$ws.Range("J89").Value = 1499.5
$ws.Range("K89").Value = 19765.385
$ws.Range("L89").Value = 7497.5
$ws.Range("M89").Value = -14149.385
$ws.Range("N89").Value = -18729.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3220.9
$ws.Range("I58").Value = 3213.625
$ws.Range("J58").Value = 3250
$ws.Range("K58").Value = 3213.625
$ws.Range("L58").Value = 3250
$ws.Range("M58").Value = -3010.625
$ws.Range("N58").Value = -3656

$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376

$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880

$ws.Range("H134").Value = 2975
$ws.Range("I134").Value = 2975
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8925
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6390
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 3220.9
$ws.Range("I136").Value = 3213.625
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 9640.875
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -7090.875
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 95.30768999999999
$ws.Range("J2").Value = 110.85714
$ws.Range("L2").Value = 665.14284
$ws.Range("N2").Value = -891.14284

$ws.Range("H4").Value = 3427135.5
$ws.Range("I4").Value = 105527.5
$ws.Range("K4").Value = 316582.5
$ws.Range("M4").Value = -316470.5

$ws.Range("H12").Value = 330.3846
$ws.Range("J12").Value = 310.6
$ws.Range("L12").Value = 931.8000000000001
$ws.Range("N12").Value = -1277.8

$ws.Range("H56").Value = 10349.553
$ws.Range("I56").Value = 10349.553
$ws.Range("K56").Value = 10349.553
$ws.Range("M56").Value = -9819.553

$ws.Range("H92").Value = 212.4
$ws.Range("I92").Value = 212.4
$ws.Range("K92").Value = 637.2
$ws.Range("M92").Value = 610.8

$ws.Range("H98").Value = 1095
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 45000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45576

$ws.Range("H81").Value = 45000
$ws.Range("J81").Value = 45000
$ws.Range("L81").Value = 45000
$ws.Range("N81").Value = -46996

$ws.Range("H84").Value = 45000
$ws.Range("J84").Value = 45000
$ws.Range("L84").Value = 135000
$ws.Range("N84").Value = -144984

$ws.Range("H102").Value = 2297.2
$ws.Range("I102").Value = 2297.2
$ws.Range("K102").Value = 2297.2
$ws.Range("M102").Value = -675.1999999999998

$ws.Range("H107").Value = 3765.5557
$ws.Range("I107").Value = 5250
$ws.Range("J107").Value = 3341.4285
$ws.Range("K107").Value = 5250
$ws.Range("L107").Value = 3341.4285
$ws.Range("M107").Value = -3330
$ws.Range("N107").Value = -7181.4285

$ws.Range("H126").Value = 2624.75
$ws.Range("I126").Value = 2624.75
$ws.Range("K126").Value = 7874.25
$ws.Range("M126").Value = -5404.25

$ws.Range("H132").Value = 2998.6667
$ws.Range("I132").Value = 2998.6667
$ws.Range("K132").Value = 8996.000100000001
$ws.Range("M132").Value = -6466.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6199.5264
$ws.Range("I7").Value = 6181.125
$ws.Range("K7").Value = 6181.125
$ws.Range("M7").Value = -6069.125

$ws.Range("H22").Value = 885.5714
$ws.Range("I22").Value = 325
$ws.Range("J22").Value = 1633
$ws.Range("K22").Value = 325
$ws.Range("L22").Value = 1633
$ws.Range("M22").Value = -30
$ws.Range("N22").Value = -2223

$ws.Range("H27").Value = 885.5714
$ws.Range("I27").Value = 325
$ws.Range("J27").Value = 1633
$ws.Range("K27").Value = 325
$ws.Range("L27").Value = 1633
$ws.Range("M27").Value = -218
$ws.Range("N27").Value = -1847

$ws.Range("H36").Value = 56666.668
$ws.Range("J36").Value = 56666.668
$ws.Range("L36").Value = 56666.668
$ws.Range("N36").Value = -57790.668

$ws.Range("H40").Value = 3974
$ws.Range("I40").Value = 3982.3333
$ws.Range("J40").Value = 3949
$ws.Range("K40").Value = 3982.3333
$ws.Range("L40").Value = 3949
$ws.Range("M40").Value = -3846.3333
$ws.Range("N40").Value = -4221

$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 1400
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1400
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -651
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 1400
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -3256
$ws.Range("N71").ClearContents()

$ws.Range("H122").Value = 3884.3333
$ws.Range("I122").Value = 3478.2666
$ws.Range("K122").Value = 10434.7998
$ws.Range("M122").Value = -7984.799800000001

$ws.Range("H126").Value = 6199.5264
$ws.Range("I126").Value = 6181.125
$ws.Range("K126").Value = 18543.375
$ws.Range("M126").Value = -16073.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H62").Value = 12678.571
$ws.Range("I62").Value = 9292.714
$ws.Range("K62").Value = 9292.714
$ws.Range("M62").Value = -8668.714

$ws.Range("H65").Value = 12678.571
$ws.Range("I65").Value = 9292.714
$ws.Range("K65").Value = 46463.57
$ws.Range("M65").Value = -43343.57

$ws.Range("H122").Value = 3441.762
$ws.Range("I122").Value = 3282.1667
$ws.Range("K122").Value = 9846.500100000001
$ws.Range("M122").Value = -7396.500100000001

$ws.Range("H132").Value = 4625.3887
$ws.Range("I132").Value = 4518.357
$ws.Range("K132").Value = 13555.071
$ws.Range("M132").Value = -11025.071

$ws.Range("H135").Value = 164857.5
$ws.Range("J135").Value = 164857.5
$ws.Range("L135").Value = 164857.5
$ws.Range("N135").Value = -174997.5
